$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 106, pushing existing rows 106-117 down to 107-118.
$ws.Rows.Item(106).Insert()

# Fill in the new row 106 with the same constant values as its neighbours,
# plus the new data point for this record.
$ws.Cells.Item(106, 1).Value = 3
$ws.Cells.Item(106, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(106, 3).Value = "Coquimbo"
$ws.Cells.Item(106, 4).Value = Get-Date -Year 2021 -Month 9 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(106, 5).Value = 5
$ws.Cells.Item(106, 6).Value = 100112010
$ws.Cells.Item(106, 7).Value = "Achicoria"
$ws.Cells.Item(106, 8).Value = "Sin especificar"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 115
$ws.Cells.Item(106, 11).Value = 6500
$ws.Cells.Item(106, 12).Value = 6800
$ws.Cells.Item(106, 13).Value = 6657
$ws.Cells.Item(106, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(106, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(106, 16).Value = 416
$ws.Cells.Item(106, 17).Value = 16
$ws.Cells.Item(106, 18).Value = "Hortaliza"
